$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.294.42"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.69"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.07"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0690"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.049.90"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.796.01"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.04"
$ws.Range("E14").Value = "  -3.13%  "

$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.271.95"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.20"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.08"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("E19").Value = "  +2.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.58"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.42"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.36"
$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("E32").Value = "  +2.39%  "

$ws.Range("E33").Value = "  +6.38%  "

$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.442.99"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  +8.30%  "

$ws.Range("E37").Value = "  +1.88%  "

$ws.Range("E38").Value = "  +1.61%  "

$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.34"
$ws.Range("E40").Value = "  +2.05%  "

$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.12"
$ws.Range("E42").Value = "  +6.32%  "

$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("E45").Value = "  +1.92%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.945.08"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.48"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("E51").Value = "  -6.26%  "
